$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$aVals = @(46056,46056.01041666666,46056.02083333334,46056.03125,46056.04166666666,46056.05208333334,46056.0625,46056.07291666666,46056.08333333334,46056.09375,46056.10416666666,46056.11458333334,46056.125,46056.13541666666,46056.14583333334,46056.15625,46056.16666666666,46056.17708333334,46056.1875,46056.19791666666,46056.20833333334,46056.21875,46056.22916666666,46056.23958333334,46056.25,46056.26041666666,46056.27083333334,46056.28125,46056.29166666666,46056.30208333334,46056.3125,46056.32291666666,46056.33333333334,46056.34375,46056.35416666666,46056.36458333334,46056.375,46056.38541666666,46056.39583333334,46056.40625,46056.41666666666,46056.42708333334,46056.4375,46056.44791666666,46056.45833333334,46056.46875,46056.47916666666,46056.48958333334,46056.5,46056.51041666666,46056.52083333334,46056.53125,46056.54166666666,46056.55208333334,46056.5625,46056.57291666666,46056.58333333334,46056.59375,46056.60416666666,46056.61458333334,46056.625,46056.63541666666,46056.64583333334,46056.65625,46056.66666666666,46056.67708333334,46056.6875,46056.69791666666,46056.70833333334,46056.71875,46056.72916666666,46056.73958333334,46056.75,46056.76041666666,46056.77083333334,46056.78125,46056.79166666666,46056.80208333334,46056.8125,46056.82291666666,46056.83333333334,46056.84375,46056.85416666666,46056.86458333334,46056.875,46056.88541666666,46056.89583333334,46056.90625,46056.91666666666,46056.92708333334,46056.9375,46056.94791666666,46056.95833333334,46056.96875,46056.97916666666,46056.98958333334)
$bVals = @(6460,6420,6410,6390,6360,6340,6330,6330,6330,6330,6330,6340,6360,6370,6380,6410,6450,6530,6640,6790,6980,7170,7370,7570,7820,8020,8200,8350,8500,8610,8680,8720,8720,8690,8640,8560,8460,8350,8250,8130,8020,7930,7870,7800,7740,7710,7690,7680,7690,7700,7700,7710,7710,7730,7760,7800,7840,7880,7920,7970,8040,8110,8190,8290,8400,8500,8570,8670,8760,8800,8820,8820,8780,8750,8730,8670,8610,8550,8470,8390,8270,8150,8020,7880,7680,7530,7360,7230,7110,7000,6890,6770,6770,6700,6650,6590)

for ($i = 0; $i -lt $aVals.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $aVals[$i]
    $ws.Cells.Item($row, 2).Value = $bVals[$i]
}
